# data archive for publication
# Re-express the collected GPS coordinates (lat/long) for the bumble bee
# voucher records on Sheet1 with their corrected/published values, drop the
# ad-hoc "black text" style that had been applied to the WI site rows
# (G5:H11) so they fall back to the workbook's Normal style, resize the
# latitude/longitude columns to a fixed width, and leave the selection where
# the editor left off before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected latitude / longitude values -------------------------------
$coords = @{
    2  = @{ G = 41.925926999999994; H = -89.349541000000002 }
    3  = @{ G = 41.923929000000001; H = -89.346551000000005 }
    4  = @{ G = 41.892828000000002; H = -89.366168999999999 }
    5  = @{ G = 43.616919000000003; H = -89.260760000000005 }
    6  = @{ G = 43.699061999999998; H = -89.396114999999995 }
    7  = @{ G = 43.699905999999999; H = -89.394874999999999 }
    8  = @{ G = 43.699953999999998; H = -89.384244000000010 }
    9  = @{ G = 43.699192000000004; H = -89.385835999999998 }
    10 = @{ G = 43.613513000000005; H = -89.261359999999996 }
    11 = @{ G = 43.614747999999999; H = -89.258881000000002 }
}

foreach ($row in $coords.Keys) {
    $ws.Range("G$row").Value = $coords[$row].G
    $ws.Range("H$row").Value = $coords[$row].H
}

# --- Drop the explicit black-font style on the WI rows (G5:H11) ----------
# These cells previously carried a one-off cell style (fontId pointing at an
# explicit black RGB color); restoring "Normal" drops that unused xf/font
# from the style table and matches the rest of the column.
$ws.Range("G5:H11").Style = "Normal"

# --- Columns G:H now get a fixed width instead of autofit ----------------
$ws.Range("G:H").ColumnWidth = 9.7

# --- Leave the selection on H8, matching the last-saved cursor position --
[void]$ws.Range("H8").Select()
